$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.445647641019636
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 35.65327920106175
$ws.Range("B3").Value = 0.00000009142958989905026
$ws.Range("C3").Value = 0.0000006194867796516235
$ws.Range("D3").Value = 3.223369029078222
$ws.Range("E3").Value = 2797.565817734744
$ws.Range("G3").Value = 2800.789187474738
